# "New helper class and diffrentlogic for assignment test added"
#
# - LMSData: refreshed Libraries / Security Domains / Create Group1 test
#   data, cleared the old Assignment Profile ID values, and appended a new
#   "Names" column (S).
# - Added a new "Domain" lookup worksheet (Names / CORE TEST A / CORE TEST B
#   / EXTERNAL) right after LMSData, with a highlighted header cell.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LMSData")

# ---------------------------------------------------------------------
# LMSData: header row - add the new "Names" column header (S1), copying
# the formatting used by the other header cells.
# ---------------------------------------------------------------------
$ws1.Range("R1").Copy()
$ws1.Range("S1").PasteSpecial(-4122)
$ws1.Cells.Item(1, 19).Value = "Names"

# ---------------------------------------------------------------------
# LMSData: the "Assignment Profile ID" values (column C, rows 2-4) are no
# longer populated - clear contents but keep the bordered formatting.
# ---------------------------------------------------------------------
$ws1.Cells.Item(2, 3).ClearContents()
$ws1.Cells.Item(3, 3).ClearContents()
$ws1.Cells.Item(4, 3).ClearContents()

# ---------------------------------------------------------------------
# LMSData: refreshed Libraries (I), Security Domains (J) and Create
# Group1 (K) test values.
# ---------------------------------------------------------------------
$ws1.Cells.Item(2, 9).Value = "COM_TEST;TestLirbrary"
$ws1.Cells.Item(3, 9).Value = "20200722;GENERAL"
$ws1.Cells.Item(4, 9).Value = "Gamification_Catalog;COM_TEST"

$ws1.Cells.Item(2, 10).Value = "CORE TEST A;EXTERNAL"
$ws1.Cells.Item(3, 10).Value = "CORE TEST A;CORE TEST B"
$ws1.Cells.Item(4, 10).Value = "CORE TEST A;CORE TEST B"

$ws1.Cells.Item(2, 11).Value = "TestGroup1;TestGroup2"
$ws1.Cells.Item(3, 11).Value = "TestGroup2"
$ws1.Cells.Item(4, 11).Value = "TestGroup3"

# ---------------------------------------------------------------------
# LMSData: new "Names" column (S) values for rows 2-4, formatted the same
# way as the other data cells in that block (copy the bordered format
# from the neighbouring Libraries column).
# ---------------------------------------------------------------------
$ws1.Range("I2").Copy()
$ws1.Range("S2").PasteSpecial(-4122)
$ws1.Range("I3").Copy()
$ws1.Range("S3").PasteSpecial(-4122)
$ws1.Range("I4").Copy()
$ws1.Range("S4").PasteSpecial(-4122)

$ws1.Cells.Item(2, 19).Value = "CORE TEST A;CORE TEST B;External"
$ws1.Cells.Item(3, 19).Value = "CORE TEST B"
$ws1.Cells.Item(4, 19).Value = "EXTERNAL"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# LMSData: widen the columns whose text grew, and the new Names column.
# ---------------------------------------------------------------------
$ws1.Columns.Item(9).ColumnWidth = 29
$ws1.Columns.Item(10).ColumnWidth = 21.825
$ws1.Columns.Item(11).ColumnWidth = 20.325
$ws1.Columns.Item(19).ColumnWidth = 31

# Scroll the view over and move the active selection, matching the
# author's on-screen state when the workbook was saved.
$ws1.Range("K9").Select()

# ---------------------------------------------------------------------
# Add the new "Domain" lookup worksheet directly after LMSData.
# ---------------------------------------------------------------------
$wsDomain = $wb.Worksheets.Add($null, $ws1)
$wsDomain.Name = "Domain"

$wsDomain.Cells.Item(1, 1).Value = "Names"
$wsDomain.Cells.Item(2, 1).Value = "CORE TEST A"
$wsDomain.Cells.Item(3, 1).Value = "CORE TEST B"
$wsDomain.Cells.Item(4, 1).Value = "EXTERNAL"

# Bordered formatting for the data rows (matches the LMSData sheet look).
$ws1.Range("N2").Copy()
$wsDomain.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header cell: bordered + highlighted with a yellow fill.
$ws1.Range("D2").Copy()
$wsDomain.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsDomain.Range("A1").Interior.Color = 65535

$wsDomain.Columns.Item(1).ColumnWidth = 11.36328125
$wsDomain.Columns.Item(2).ColumnWidth = 21.825

$wsDomain.Range("A1").Select()

$ws1.Activate()
